$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-93 and add new rows 94-98 with shifted date/load data
$ws.Cells.Item(2,1).Value = 45469
$ws.Cells.Item(2,2).Value = 5730
$ws.Cells.Item(3,1).Value = 45469.01041666666
$ws.Cells.Item(3,2).Value = 5650
$ws.Cells.Item(4,1).Value = 45469.02083333334
$ws.Cells.Item(4,2).Value = 5580
$ws.Cells.Item(5,1).Value = 45469.03125
$ws.Cells.Item(5,2).Value = 5520
$ws.Cells.Item(6,1).Value = 45469.04166666666
$ws.Cells.Item(6,2).Value = 5480
$ws.Cells.Item(7,1).Value = 45469.05208333334
$ws.Cells.Item(7,2).Value = 5440
$ws.Cells.Item(8,1).Value = 45469.0625
$ws.Cells.Item(8,2).Value = 5420
$ws.Cells.Item(9,1).Value = 45469.07291666666
$ws.Cells.Item(9,2).Value = 5410
$ws.Cells.Item(10,1).Value = 45469.08333333334
$ws.Cells.Item(10,2).Value = 5390
$ws.Cells.Item(11,1).Value = 45469.09375
$ws.Cells.Item(11,2).Value = 5380
$ws.Cells.Item(12,1).Value = 45469.10416666666
$ws.Cells.Item(12,2).Value = 5370
$ws.Cells.Item(13,1).Value = 45469.11458333334
$ws.Cells.Item(13,2).Value = 5360
$ws.Cells.Item(14,1).Value = 45469.125
$ws.Cells.Item(14,2).Value = 5350
$ws.Cells.Item(15,1).Value = 45469.13541666666
$ws.Cells.Item(15,2).Value = 5350
$ws.Cells.Item(16,1).Value = 45469.14583333334
$ws.Cells.Item(16,2).Value = 5360
$ws.Cells.Item(17,1).Value = 45469.15625
$ws.Cells.Item(17,2).Value = 5370
$ws.Cells.Item(18,1).Value = 45469.16666666666
$ws.Cells.Item(18,2).Value = 5380
$ws.Cells.Item(19,1).Value = 45469.17708333334
$ws.Cells.Item(19,2).Value = 5390
$ws.Cells.Item(20,1).Value = 45469.1875
$ws.Cells.Item(20,2).Value = 5420
$ws.Cells.Item(21,1).Value = 45469.19791666666
$ws.Cells.Item(21,2).Value = 5480
$ws.Cells.Item(22,1).Value = 45469.20833333334
$ws.Cells.Item(22,2).Value = 5580
$ws.Cells.Item(23,1).Value = 45469.21875
$ws.Cells.Item(23,2).Value = 5710
$ws.Cells.Item(24,1).Value = 45469.22916666666
$ws.Cells.Item(24,2).Value = 5850
$ws.Cells.Item(25,1).Value = 45469.23958333334
$ws.Cells.Item(25,2).Value = 6000
$ws.Cells.Item(26,1).Value = 45469.25
$ws.Cells.Item(26,2).Value = 6150
$ws.Cells.Item(27,1).Value = 45469.26041666666
$ws.Cells.Item(27,2).Value = 6280
$ws.Cells.Item(28,1).Value = 45469.27083333334
$ws.Cells.Item(28,2).Value = 6400
$ws.Cells.Item(29,1).Value = 45469.28125
$ws.Cells.Item(29,2).Value = 6500
$ws.Cells.Item(30,1).Value = 45469.29166666666
$ws.Cells.Item(30,2).Value = 6580
$ws.Cells.Item(31,1).Value = 45469.30208333334
$ws.Cells.Item(31,2).Value = 6620
$ws.Cells.Item(32,1).Value = 45469.3125
$ws.Cells.Item(32,2).Value = 6640
$ws.Cells.Item(33,1).Value = 45469.32291666666
$ws.Cells.Item(33,2).Value = 6640
$ws.Cells.Item(34,1).Value = 45469.33333333334
$ws.Cells.Item(34,2).Value = 6620
$ws.Cells.Item(35,1).Value = 45469.34375
$ws.Cells.Item(35,2).Value = 6590
$ws.Cells.Item(36,1).Value = 45469.35416666666
$ws.Cells.Item(36,2).Value = 6560
$ws.Cells.Item(37,1).Value = 45469.36458333334
$ws.Cells.Item(37,2).Value = 6530
$ws.Cells.Item(38,1).Value = 45469.375
$ws.Cells.Item(38,2).Value = 6490
$ws.Cells.Item(39,1).Value = 45469.38541666666
$ws.Cells.Item(39,2).Value = 6470
$ws.Cells.Item(40,1).Value = 45469.39583333334
$ws.Cells.Item(40,2).Value = 6440
$ws.Cells.Item(41,1).Value = 45469.40625
$ws.Cells.Item(41,2).Value = 6420
$ws.Cells.Item(42,1).Value = 45469.41666666666
$ws.Cells.Item(42,2).Value = 6390
$ws.Cells.Item(43,1).Value = 45469.42708333334
$ws.Cells.Item(43,2).Value = 6370
$ws.Cells.Item(44,1).Value = 45469.4375
$ws.Cells.Item(44,2).Value = 6370
$ws.Cells.Item(45,1).Value = 45469.44791666666
$ws.Cells.Item(45,2).Value = 6380
$ws.Cells.Item(46,1).Value = 45469.45833333334
$ws.Cells.Item(46,2).Value = 6420
$ws.Cells.Item(47,1).Value = 45469.46875
$ws.Cells.Item(47,2).Value = 6450
$ws.Cells.Item(48,1).Value = 45469.47916666666
$ws.Cells.Item(48,2).Value = 6490
$ws.Cells.Item(49,1).Value = 45469.48958333334
$ws.Cells.Item(49,2).Value = 6520
$ws.Cells.Item(50,1).Value = 45469.5
$ws.Cells.Item(50,2).Value = 6550
$ws.Cells.Item(51,1).Value = 45469.51041666666
$ws.Cells.Item(51,2).Value = 6570
$ws.Cells.Item(52,1).Value = 45469.52083333334
$ws.Cells.Item(52,2).Value = 6590
$ws.Cells.Item(53,1).Value = 45469.53125
$ws.Cells.Item(53,2).Value = 6610
$ws.Cells.Item(54,1).Value = 45469.54166666666
$ws.Cells.Item(54,2).Value = 6620
$ws.Cells.Item(55,1).Value = 45469.55208333334
$ws.Cells.Item(55,2).Value = 6620
$ws.Cells.Item(56,1).Value = 45469.5625
$ws.Cells.Item(56,2).Value = 6630
$ws.Cells.Item(57,1).Value = 45469.57291666666
$ws.Cells.Item(57,2).Value = 6640
$ws.Cells.Item(58,1).Value = 45469.58333333334
$ws.Cells.Item(58,2).Value = 6650
$ws.Cells.Item(59,1).Value = 45469.59375
$ws.Cells.Item(59,2).Value = 6660
$ws.Cells.Item(60,1).Value = 45469.60416666666
$ws.Cells.Item(60,2).Value = 6670
$ws.Cells.Item(61,1).Value = 45469.61458333334
$ws.Cells.Item(61,2).Value = 6680
$ws.Cells.Item(62,1).Value = 45469.625
$ws.Cells.Item(62,2).Value = 6730
$ws.Cells.Item(63,1).Value = 45469.63541666666
$ws.Cells.Item(63,2).Value = 6770
$ws.Cells.Item(64,1).Value = 45469.64583333334
$ws.Cells.Item(64,2).Value = 6820
$ws.Cells.Item(65,1).Value = 45469.65625
$ws.Cells.Item(65,2).Value = 6870
$ws.Cells.Item(66,1).Value = 45469.66666666666
$ws.Cells.Item(66,2).Value = 6920
$ws.Cells.Item(67,1).Value = 45469.67708333334
$ws.Cells.Item(67,2).Value = 6960
$ws.Cells.Item(68,1).Value = 45469.6875
$ws.Cells.Item(68,2).Value = 7010
$ws.Cells.Item(69,1).Value = 45469.69791666666
$ws.Cells.Item(69,2).Value = 7050
$ws.Cells.Item(70,1).Value = 45469.70833333334
$ws.Cells.Item(70,2).Value = 7110
$ws.Cells.Item(71,1).Value = 45469.71875
$ws.Cells.Item(71,2).Value = 7160
$ws.Cells.Item(72,1).Value = 45469.72916666666
$ws.Cells.Item(72,2).Value = 7220
$ws.Cells.Item(73,1).Value = 45469.73958333334
$ws.Cells.Item(73,2).Value = 7290
$ws.Cells.Item(74,1).Value = 45469.75
$ws.Cells.Item(74,2).Value = 7360
$ws.Cells.Item(75,1).Value = 45469.76041666666
$ws.Cells.Item(75,2).Value = 7430
$ws.Cells.Item(76,1).Value = 45469.77083333334
$ws.Cells.Item(76,2).Value = 7470
$ws.Cells.Item(77,1).Value = 45469.78125
$ws.Cells.Item(77,2).Value = 7510
$ws.Cells.Item(78,1).Value = 45469.79166666666
$ws.Cells.Item(78,2).Value = 7540
$ws.Cells.Item(79,1).Value = 45469.80208333334
$ws.Cells.Item(79,2).Value = 7570
$ws.Cells.Item(80,1).Value = 45469.8125
$ws.Cells.Item(80,2).Value = 7590
$ws.Cells.Item(81,1).Value = 45469.82291666666
$ws.Cells.Item(81,2).Value = 7620
$ws.Cells.Item(82,1).Value = 45469.83333333334
$ws.Cells.Item(82,2).Value = 7640
$ws.Cells.Item(83,1).Value = 45469.84375
$ws.Cells.Item(83,2).Value = 7650
$ws.Cells.Item(84,1).Value = 45469.85416666666
$ws.Cells.Item(84,2).Value = 7630
$ws.Cells.Item(85,1).Value = 45469.86458333334
$ws.Cells.Item(85,2).Value = 7570
$ws.Cells.Item(86,1).Value = 45469.875
$ws.Cells.Item(86,2).Value = 7480
$ws.Cells.Item(87,1).Value = 45469.88541666666
$ws.Cells.Item(87,2).Value = 7350
$ws.Cells.Item(88,1).Value = 45469.89583333334
$ws.Cells.Item(88,2).Value = 7200
$ws.Cells.Item(89,1).Value = 45469.90625
$ws.Cells.Item(89,2).Value = 7020
$ws.Cells.Item(90,1).Value = 45469.91666666666
$ws.Cells.Item(90,2).Value = 6820
$ws.Cells.Item(91,1).Value = 45469.92708333334
$ws.Cells.Item(91,2).Value = 6660
$ws.Cells.Item(92,1).Value = 45469.9375
$ws.Cells.Item(92,2).Value = 6550
$ws.Cells.Item(93,1).Value = 45469.94791666666
$ws.Cells.Item(93,2).Value = 6440
$ws.Cells.Item(94,1).Value = 45469.95833333334
$ws.Cells.Item(94,2).Value = 5970
$ws.Cells.Item(95,1).Value = 45469.96875
$ws.Cells.Item(95,2).Value = 5900
$ws.Cells.Item(96,1).Value = 45469.97916666666
$ws.Cells.Item(96,2).Value = 5840
$ws.Cells.Item(97,1).Value = 45469.98958333334
$ws.Cells.Item(97,2).Value = 5780
$ws.Cells.Item(98,1).Value = 45470
$ws.Cells.Item(98,2).Value = 5710

# Ensure new rows (94-98) use the same date/time number format as column A elsewhere
$ws.Range("A94:A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
